$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.138059000000001
$ws.Range("H2").Value = 18.414177
$ws.Range("I2").Value = 0.08535364925338249
$ws.Range("J2").Value = 0.08535364925338247
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 28.72417333333333
$ws.Range("N2").Value = 86.17251999999999
$ws.Range("O2").Value = 0.4233259107972328
$ws.Range("P2").Value = 0.4233259107972328
$ws.Range("Q2").Value = 176.3106706462267
$ws.Range("R2").Value = 1586.79603581604
$ws.Range("S2").Value = 0.03613241131005569
$ws.Range("T2").Value = 0.03613241131005569

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.138059000000001
$ws.Range("H3").Value = 18.414177
$ws.Range("I3").Value = 0.08535364925338249
$ws.Range("J3").Value = 0.08535364925338247
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 30.56986233333333
$ws.Range("N3").Value = 91.709587
$ws.Range("O3").Value = 0.4505269713084062
$ws.Range("P3").Value = 0.4505269713084062
$ws.Range("Q3").Value = 187.6396186238777
$ws.Range("R3").Value = 1688.756567614899
$ws.Range("S3").Value = 0.03845412108824642
$ws.Range("T3").Value = 0.03845412108824642

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.138059000000001
$ws.Range("H4").Value = 18.414177
$ws.Range("I4").Value = 0.08535364925338249
$ws.Range("J4").Value = 0.08535364925338247
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.559531999999999
$ws.Range("N4").Value = 25.678596
$ws.Range("O4").Value = 0.126147117894361
$ws.Range("P4").Value = 0.126147117894361
$ws.Range("Q4").Value = 52.538912428388
$ws.Range("R4").Value = 472.850211855492
$ws.Range("S4").Value = 0.01076711685508038
$ws.Range("T4").Value = 0.01076711685508037

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.15353
$ws.Range("H5").Value = 111.46059
$ws.Range("I5").Value = 0.5166436764692264
$ws.Range("J5").Value = 0.5166436764692264
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 28.72417333333333
$ws.Range("N5").Value = 86.17251999999999
$ws.Range("O5").Value = 0.4233259107972328
$ws.Range("P5").Value = 0.4233259107972328
$ws.Range("Q5").Value = 1067.2044356652
$ws.Range("R5").Value = 9604.839920986798
$ws.Range("S5").Value = 0.2187086548989661
$ws.Range("T5").Value = 0.2187086548989661

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.15353
$ws.Range("H6").Value = 111.46059
$ws.Range("I6").Value = 0.5166436764692264
$ws.Range("J6").Value = 0.5166436764692264
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.56986233333333
$ws.Range("N6").Value = 91.709587
$ws.Range("O6").Value = 0.4505269713084062
$ws.Range("P6").Value = 0.4505269713084062
$ws.Range("Q6").Value = 1135.77829729737
$ws.Range("R6").Value = 10222.00467567633
$ws.Range("S6").Value = 0.2327619108053207
$ws.Range("T6").Value = 0.2327619108053207

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.15353
$ws.Range("H7").Value = 111.46059
$ws.Range("I7").Value = 0.5166436764692264
$ws.Range("J7").Value = 0.5166436764692264
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.559531999999999
$ws.Range("N7").Value = 25.678596
$ws.Range("O7").Value = 0.126147117894361
$ws.Range("P7").Value = 0.126147117894361
$ws.Range("Q7").Value = 318.0168289479599
$ws.Range("R7").Value = 2862.15146053164
$ws.Range("S7").Value = 0.06517311076493959
$ws.Range("T7").Value = 0.06517311076493959

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.621669
$ws.Range("H8").Value = 85.86500699999999
$ws.Range("I8").Value = 0.3980026742773913
$ws.Range("J8").Value = 0.3980026742773912
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 28.72417333333333
$ws.Range("N8").Value = 86.17251999999999
$ws.Range("O8").Value = 0.4233259107972328
$ws.Range("P8").Value = 0.4233259107972328
$ws.Range("Q8").Value = 822.1337814452932
$ws.Range("R8").Value = 7399.204033007639
$ws.Range("S8").Value = 0.168484844588211
$ws.Range("T8").Value = 0.168484844588211

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.621669
$ws.Range("H9").Value = 85.86500699999999
$ws.Range("I9").Value = 0.3980026742773913
$ws.Range("J9").Value = 0.3980026742773912
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.56986233333333
$ws.Range("N9").Value = 91.709587
$ws.Range("O9").Value = 0.4505269713084062
$ws.Range("P9").Value = 0.4505269713084062
$ws.Range("Q9").Value = 874.9604810802342
$ws.Range("R9").Value = 7874.644329722109
$ws.Range("S9").Value = 0.1793109394148392
$ws.Range("T9").Value = 0.1793109394148392

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.621669
$ws.Range("H10").Value = 85.86500699999999
$ws.Range("I10").Value = 0.3980026742773913
$ws.Range("J10").Value = 0.3980026742773912
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.559531999999999
$ws.Range("N10").Value = 25.678596
$ws.Range("O10").Value = 0.126147117894361
$ws.Range("P10").Value = 0.126147117894361
$ws.Range("Q10").Value = 244.9880916989079
$ws.Range("R10").Value = 2204.892825290171
$ws.Range("S10").Value = 0.05020689027434103
$ws.Range("T10").Value = 0.05020689027434102
